# Auto-generated Excel COM-interop script
# Commit: update scripts wuth new tpm
# Recomputed NATMI ligand-receptor edge-weight table (Ncam1 -> Fgfr1) values
# after the underlying per-cluster TPM expression inputs were regenerated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("G2", "0.3912683333333333")
    ,@("H2", "1.173805")
    ,@("I2", "0.004923718964983145")
    ,@("J2", "0.004923718964983145")
    ,@("M2", "7.955277333333332")
    ,@("N2", "23.865832")
    ,@("O2", "0.05015625076675284")
    ,@("P2", "0.05015625076675283")
    ,@("Q2", "3.112648103417777")
    ,@("R2", "28.01383293076")
    ,@("S2", "0.0002469552831127114")
    ,@("T2", "0.0002469552831127114")
    ,@("G3", "0.3912683333333333")
    ,@("H3", "1.173805")
    ,@("I3", "0.004923718964983145")
    ,@("J3", "0.004923718964983145")
    ,@("M3", "82.48060333333333")
    ,@("O3", "0.520021823355633")
    ,@("P3", "0.520021823355633")
    ,@("Q3", "32.27204819856111")
    ,@("R3", "290.44843378705")
    ,@("S3", "0.002560441313861245")
    ,@("T3", "0.002560441313861245")
    ,@("G4", "0.3912683333333333")
    ,@("H4", "1.173805")
    ,@("I4", "0.004923718964983145")
    ,@("J4", "0.004923718964983145")
    ,@("M4", "66.90297433333333")
    ,@("N4", "200.708923")
    ,@("O4", "0.4218083439585467")
    ,@("P4", "0.4218083439585465")
    ,@("Q4", "26.17701526244611")
    ,@("R4", "235.593137362015")
    ,@("S4", "0.00207686574273683")
    ,@("T4", "0.002076865742736829")
    ,@("G5", "0.3912683333333333")
    ,@("H5", "1.173805")
    ,@("I5", "0.004923718964983145")
    ,@("J5", "0.004923718964983145")
    ,@("M5", "1.271033333333333")
    ,@("N5", "3.8131")
    ,@("O5", "0.008013581919067616")
    ,@("P5", "0.008013581919067614")
    ,@("Q5", "0.4973150939444445")
    ,@("R5", "4.475835845500001")
    ,@("S5", "3.945662527235925e-05")
    ,@("T5", "3.945662527235924e-05")
    ,@("I6", "0.07888477275715973")
    ,@("J6", "0.07888477275715973")
    ,@("M6", "7.955277333333332")
    ,@("N6", "23.865832")
    ,@("O6", "0.05015625076675284")
    ,@("P6", "0.05015625076675283")
    ,@("Q6", "49.86891820133688")
    ,@("R6", "448.820263812032")
    ,@("S6", "0.003956564444086416")
    ,@("T6", "0.003956564444086416")
    ,@("I7", "0.07888477275715973")
    ,@("J7", "0.07888477275715973")
    ,@("M7", "82.48060333333333")
    ,@("O7", "0.520021823355633")
    ,@("P7", "0.520021823355633")
    ,@("R7", "4653.384740256561")
    ,@("S7", "0.04102180336417297")
    ,@("T7", "0.04102180336417297")
    ,@("I8", "0.07888477275715973")
    ,@("J8", "0.07888477275715973")
    ,@("M8", "66.90297433333333")
    ,@("N8", "200.708923")
    ,@("O8", "0.4218083439585467")
    ,@("P8", "0.4218083439585465")
    ,@("Q8", "419.3919098804276")
    ,@("R8", "3774.527188923848")
    ,@("S8", "0.03327425536024382")
    ,@("T8", "0.03327425536024382")
    ,@("I9", "0.07888477275715973")
    ,@("J9", "0.07888477275715973")
    ,@("M9", "1.271033333333333")
    ,@("N9", "3.8131")
    ,@("O9", "0.008013581919067616")
    ,@("P9", "0.008013581919067614")
    ,@("Q9", "7.967674120622224")
    ,@("R9", "71.70906708560001")
    ,@("S9", "0.0006321495886565328")
    ,@("T9", "0.0006321495886565327")
    ,@("G10", "72.68848166666666")
    ,@("H10", "218.065445")
    ,@("I10", "0.9147115297293749")
    ,@("J10", "0.9147115297293749")
    ,@("M10", "7.955277333333332")
    ,@("N10", "23.865832")
    ,@("O10", "0.05015625076675284")
    ,@("P10", "0.05015625076675283")
    ,@("Q10", "578.2570305972488")
    ,@("R10", "5204.313275375239")
    ,@("S10", "0.04587850086434662")
    ,@("T10", "0.04587850086434662")
    ,@("G11", "72.68848166666666")
    ,@("H11", "218.065445")
    ,@("I11", "0.9147115297293749")
    ,@("J11", "0.9147115297293749")
    ,@("M11", "82.48060333333333")
    ,@("O11", "0.520021823355633")
    ,@("P11", "0.520021823355633")
    ,@("Q11", "5995.389823250605")
    ,@("R11", "53958.50840925545")
    ,@("S11", "0.4756699575342898")
    ,@("T11", "0.4756699575342898")
    ,@("G12", "72.68848166666666")
    ,@("H12", "218.065445")
    ,@("I12", "0.9147115297293749")
    ,@("J12", "0.9147115297293749")
    ,@("M12", "66.90297433333333")
    ,@("N12", "200.708923")
    ,@("O12", "0.4218083439585467")
    ,@("P12", "0.4218083439585465")
    ,@("Q12", "4863.07562327397")
    ,@("R12", "43767.68060946573")
    ,@("S12", "0.3858329555549365")
    ,@("T12", "0.3858329555549364")
    ,@("G13", "72.68848166666666")
    ,@("H13", "218.065445")
    ,@("I13", "0.9147115297293749")
    ,@("J13", "0.9147115297293749")
    ,@("M13", "1.271033333333333")
    ,@("N13", "3.8131")
    ,@("O13", "0.008013581919067616")
    ,@("P13", "0.008013581919067614")
    ,@("Q13", "92.38948314772223")
    ,@("R13", "831.5053483295001")
    ,@("S13", "0.007330115775801998")
    ,@("T13", "0.007330115775801997")
    ,@("E14", "2.0")
    ,@("F14", "0.6666666666666666")
    ,@("G14", "0.117608")
    ,@("H14", "0.352824")
    ,@("I14", "0.001479978548482255")
    ,@("J14", "0.001479978548482255")
    ,@("M14", "7.955277333333332")
    ,@("N14", "23.865832")
    ,@("O14", "0.05015625076675284")
    ,@("P14", "0.05015625076675283")
    ,@("Q14", "0.9356042566186665")
    ,@("R14", "8.420438309568")
    ,@("S14", "7.423017520709086e-05")
    ,@("T14", "7.423017520709085e-05")
    ,@("E15", "2.0")
    ,@("F15", "0.6666666666666666")
    ,@("G15", "0.117608")
    ,@("H15", "0.352824")
    ,@("I15", "0.001479978548482255")
    ,@("J15", "0.001479978548482255")
    ,@("M15", "82.48060333333333")
    ,@("O15", "0.520021823355633")
    ,@("P15", "0.520021823355633")
    ,@("Q15", "9.700378796826667")
    ,@("R15", "87.30340917144001")
    ,@("S15", "0.0007696211433089652")
    ,@("T15", "0.0007696211433089652")
    ,@("E16", "2.0")
    ,@("F16", "0.6666666666666666")
    ,@("G16", "0.117608")
    ,@("H16", "0.352824")
    ,@("I16", "0.001479978548482255")
    ,@("J16", "0.001479978548482255")
    ,@("M16", "66.90297433333333")
    ,@("N16", "200.708923")
    ,@("O16", "0.4218083439585467")
    ,@("P16", "0.4218083439585465")
    ,@("Q16", "7.868325005394667")
    ,@("R16", "70.814925048552")
    ,@("S16", "0.0006242673006294736")
    ,@("T16", "0.0006242673006294734")
    ,@("E17", "2.0")
    ,@("F17", "0.6666666666666666")
    ,@("G17", "0.117608")
    ,@("H17", "0.352824")
    ,@("I17", "0.001479978548482255")
    ,@("J17", "0.001479978548482255")
    ,@("M17", "1.271033333333333")
    ,@("N17", "3.8131")
    ,@("O17", "0.008013581919067616")
    ,@("P17", "0.008013581919067614")
    ,@("Q17", "0.1494836882666667")
    ,@("R17", "1.3453531944")
    ,@("S17", "1.185992933672533e-05")
    ,@("T17", "1.185992933672533e-05")
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

Write-Host "Updated $($updates.Count) cells"
